$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.730.43"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "'1.897.38"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'246.73"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.4725"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "'0.2928"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.06515"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'22.59"
$ws.Range("E10").Value = "  +2.32%  "
$ws.Range("D11").Value = "'0.07784"
$ws.Range("D12").Value = "'1.892.86"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "'0.7406"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "'96.68"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "'5.223"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "'285.60"
$ws.Range("E16").Value = "  +3.92%  "
$ws.Range("D17").Value = "'30.732.52"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "'13.19"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").Value = "'0.000007524"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'2.133.43"
$ws.Range("D22").Value = "'5.305"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'6.270"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "'9.202"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("D26").Value = "'164.59"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").Value = "'19.00"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "'1.920"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.342"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.09783"
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("D31").Value = "'1.489"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").Value = "'4.319"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "'4.165"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "'0.04905"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("D36").Value = "'0.6981"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'2.719"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "'0.01908"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("D39").Value = "'2.851"
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("D40").Value = "'76.01"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("D41").Value = "'6.298"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("D42").Value = "'2.007"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").Value = "'0.4283"
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "'0.8296"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").Value = "'101.64"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "'9.589"
$ws.Range("E47").Value = "  +2.57%  "
$ws.Range("D48").Value = "'7.000"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").Value = "'35.40"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").Value = "'909.46"
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("D51").Value = "'0.05772"
$ws.Range("E51").Value = "  +2.19%  "
